$wb = $excel.ActiveWorkbook

# The "queries" sheet holds a query definition table. Column G ("selectionArgs")
# for the select_linked_data queries used the raw JS helper
# encodeURIComponent(...) to build the household_id URI fragment. Replace it
# with the ODK Survey helper opendatakit.encodeURIDataElement(...), which is
# the more appropriate way to build this string for use inside linked_table
# queries.
$ws = $wb.Worksheets.Item("queries")
$ws.Activate()

$oldValue = "'household_id='+encodeURIComponent(data('household_id'))"
# Leading apostrophe is a special "store as text" prefix in Excel: assigning
# a string starting with a single "'" via .Value causes Excel to swallow that
# leading character. Double it so the literal leading apostrophe survives.
$newValue = "''household_id='+opendatakit.encodeURIDataElement('household_id')"

$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    for ($c = 1; $c -le $used.Columns.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.Value2 -eq $oldValue) {
            $cell.Value = $newValue
        }
    }
}

$lastCell = $ws.Cells.Item(4, 7)
$lastCell.Select()
